$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.741.91'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '2.354.13'
$ws.Range('E3').Value = '  -4.39%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '539.83'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '137.59'
$ws.Range('E6').Value = '  -6.04%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.522'
$ws.Range('E8').Value = '  -11.12%  '
$ws.Range('D9').Value = '2.351.51'
$ws.Range('E9').Value = '  -4.42%  '
$ws.Range('E10').Value = '  -1.84%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  -3.66%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.340'
$ws.Range('E13').Value = '  -3.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '24.76'
$ws.Range('E14').Value = '  -5.19%  '
$ws.Range('D15').Value = '2.779.02'
$ws.Range('E15').Value = '  -4.35%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000162'
$ws.Range('E16').Value = '  -2.82%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '60.648.88'
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('D18').Value = '2.355.95'
$ws.Range('E18').Value = '  -3.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.57'
$ws.Range('E19').Value = '  -4.95%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '315.10'
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('E22').Value = '  -6.97%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('E24').Value = '  +2.80%  '
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.47'
$ws.Range('E26').Value = '  +11.12%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '2.475.64'
$ws.Range('E28').Value = '  -3.98%  '
$ws.Range('E29').Value = '  -7.23%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.96'
$ws.Range('E30').Value = '  -3.45%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '503.62'
$ws.Range('E31').Value = '  -8.49%  '
$ws.Range('E32').Value = '  -5.48%  '
$ws.Range('E33').Value = '  -1.79%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').Value = '  -6.14%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('E37').Value = '  -5.68%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.372'
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.26'
$ws.Range('E40').Value = '  -10.39%  '
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '138.58'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '40.10'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('E45').Value = '  -8.13%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '138.55'
$ws.Range('E46').Value = '  -5.47%  '
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('E48').Value = '  -4.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '19.47'
$ws.Range('E49').Value = '  -9.10%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.567'
$ws.Range('E50').Value = '  -3.51%  '
$ws.Range('E51').Value = '  -4.38%  '
